$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.375.17'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '1.606.09'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0850'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '1.830.95'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.608.93'
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("D17").Value = '26.377.12'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.56%  '
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("E20").Value = '  +3.52%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.07%  '
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0494'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").Value = '1.495.61'
$ws.Range("E32").Value = '  +5.30%  '
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.561'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.35%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.819'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("E43").Value = '  -4.14%  '
$ws.Range("D44").Value = '1.744.55'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.759'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0500'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0958'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  +0.25%  '
